# Generate Report for Handback
# Updates the localization-status report: marks zh-cn/de-de rows as handed
# back (in sync with en-US), refreshes the "Latest Handback DateTime"
# timestamps, clears the stale "Error Detail" warning, and widens the
# Status / Error Detail columns to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns -------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# --- zh-cn handback detail sheet ------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-13 11:01:06"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# --- de-de handback detail sheet ------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-13 11:01:15"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
